$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the "Responsable" values between the two code tasks (B26 <-> B28)
# B26: "Codage fonctions motrices" -> was Quentin, now Souf
# B28: "Codage fonctions sensorielles" -> was Souf, now Quentin
$tmp = $ws.Range("B26").Value2
$ws.Range("B26").Value2 = $ws.Range("B28").Value2
$ws.Range("B28").Value2 = $tmp

# Update the active selection to match the saved view state (B29)
$ws.Range("B29").Select()
